$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 value
$ws.Range("B3").Value = 3483

# Update D3 formula - unshare it (make it its own ordinary formula)
$ws.Range("D3").Formula = "=+D2+C3"

# Update D4 formula - becomes the new shared formula master
$ws.Range("D4").Formula = "=+D3+C4"

# Update selection to D4
$ws.Range("D4").Select()
